$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Datos"
[void]$ws.Range("E26").Select()
